$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 25, shifting existing rows 25-71 down to 26-72.
$ws.Rows("25:25").Insert()

# Populate the new row 25 with a new weekly price observation.
# Most fields mirror the (now shifted-down) row that used to be row 25,
# except for the date (D) and volume (J), which are new.
$ws.Range("A25").Value = 10
$ws.Range("B25").Value = "Vega Modelo de Temuco"
$ws.Range("C25").Value = "La Araucanía"
$ws.Range("D25").Value = 45012
$ws.Range("E25").Value = 9
$ws.Range("F25").Value = 100112010
$ws.Range("G25").Value = "Achicoria"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 10000
$ws.Range("N25").Value = "$/caja 18 unidades"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 556
$ws.Range("Q25").Value = 18
$ws.Range("R25").Value = "Hortaliza"
